# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns
# per the scraped source data for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'43.578.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'2.275.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.76%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5: Solana
$ws.Range("D5").Value = "'123.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.78%  "

# Row 6: BNB
$ws.Range("D6").Value = "'266.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "

# Row 7: XRP
$ws.Range("E7").Value = "  +2.21%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.20%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.24%  "

# Row 10: Avalanche
$ws.Range("D10").Value = "'48.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.65%  "

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.0947"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "

# Row 12: Polkadot
$ws.Range("D12").Value = "'9.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.86%  "

# Row 13: TRON
$ws.Range("E13").Value = "  -0.82%  "

# Row 14: Chainlink
$ws.Range("D14").Value = "'15.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.04%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.908"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.97%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "'2.617.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "'2.267.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.30%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "'43.534.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "

# Row 19: ShibaInu
$ws.Range("E19").Value = "  +0.66%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "'6.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "

# Row 21: Litecoin
$ws.Range("D21").Value = "'72.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.28%  "

# Row 22: ImmutableX
$ws.Range("D22").Value = "'2.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "

# Row 23: BitcoinCash
$ws.Range("D23").Value = "'235.48"
$ws.Range("D23").Style = "Normal"

# Row 24: PancakeSwap
$ws.Range("E24").Value = "  -1.73%  "

# Row 25: InternetComputer(DFINITY)
$ws.Range("D25").Value = "'9.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.85%  "

# Row 26: Cosmos
$ws.Range("D26").Value = "'11.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.31%  "

# Row 27: Dai
$ws.Range("E27").Value = "  +1.65%  "

# Row 28: InjectiveProtocol
$ws.Range("D28").Value = "'42.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.78%  "

# Row 29: WEMIXToken
$ws.Range("E29").Value = "  -0.27%  "

# Row 30: Toncoin
$ws.Range("E30").Value = "  +0.56%  "

# Row 31: Monero
$ws.Range("D31").Value = "'172.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "

# Row 32: EthereumClassic
$ws.Range("E32").Value = "  +0.32%  "

# Row 33: Hedera
$ws.Range("D33").Value = "'0.0915"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.35%  "

# Row 34: Filecoin
$ws.Range("E34").Value = "  +0.28%  "

# Row 35: Stellar
$ws.Range("E35").Value = "  +1.54%  "

# Row 36: NEARProtocol
$ws.Range("E36").Value = "  +12.56%  "

# Row 37: VeChain
$ws.Range("E37").Value = "  +4.62%  "

# Row 38: RenderToken
$ws.Range("D38").Value = "'4.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.96%  "

# Row 39: Kaspa
$ws.Range("E39").Value = "  -2.30%  "

# Row 40: LidoDAOToken
$ws.Range("D40").Value = "'2.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.63%  "

# Row 41: Celestia
$ws.Range("D41").Value = "'14.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.14%  "

# Row 42: MultiversX
$ws.Range("D42").Value = "'73.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43: Algorand
$ws.Range("E43").Value = "  -1.24%  "

# Row 44: FirstDigitalUSD
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "

# Row 45: ARBITRUM
$ws.Range("E45").Value = "  -1.17%  "

# Row 46: THORChain
$ws.Range("E46").Value = "  -11.35%  "

# Row 47: ordi
$ws.Range("D47").Value = "'74.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +37.24%  "

# Row 48: FraxShare
$ws.Range("E48").Value = "  -1.99%  "

# Row 49: TrustWalletToken
$ws.Range("E49").Value = "  -0.14%  "

# Row 50: Cronos
$ws.Range("E50").Value = "  -0.20%  "

# Row 51: Aave
$ws.Range("D51").Value = "'101.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.53%  "
